$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("B5").Value = 247.349257777547
$ws.Range("C5").Value = 259.459866270917
$ws.Range("E5").Value = 242.868070200222
$ws.Range("H5").Value = 272.784403600927
$ws.Range("I5").Value = 248.420806426136
$ws.Range("J5").Value = 217.581780325461
$ws.Range("K5").Value = 224.86166955442

# Row 6 updates
$ws.Range("B6").Value = 5.43932924266909
$ws.Range("C6").Value = 5.78615006681583
$ws.Range("E6").Value = 5.46207705182127
$ws.Range("H6").Value = 5.97658731409573
$ws.Range("I6").Value = 5.55721448027376
$ws.Range("J6").Value = 4.6981207899457
$ws.Range("K6").Value = 4.60593126685029
